$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style (bold, border, centered) from an existing header cell (H1) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Data values for columns I and J (rows 2-9)
$values = @{
    2 = @(9, 9)
    3 = @(8, 8)
    4 = @(7, 7)
    5 = @(8, 9)
    6 = @(7, 7)
    7 = @(9, 9)
    8 = @(9, 9)
    9 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
